$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    Heading1 title paragraph ("Play African Simba slot for free -
#    Review"). We clone the run/formatting structure of the existing
#    bold "title" paragraph near the end of the doc (empty leading
#    run + bold run + plain run) so the OOXML shape matches exactly,
#    then rewrite its text in place.
# ------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

# Source paragraph with the empty-run + bold-run pattern we want to copy.
$boldSource = $d.Paragraphs.Item(49)
$metaPara.Range.FormattedText = $boldSource.Range.FormattedText

# Replace the bold run's text (still bold) with "Meta description".
$metaPara.Range.Find.Execute("Play African Simba slot for free - Review", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Meta description", 2) | Out-Null

# Append the (non-bold) rest of the sentence right before the paragraph mark.
$metaEnd = $metaPara.Range.End
$insertPoint = $d.Range($metaEnd - 1, $metaEnd - 1)
$insertPoint.InsertAfter(": Read our review of African Simba, an online slot game in the African savanna theme. Play for free and win big with bonus rounds and the Gamble feature.")

# ------------------------------------------------------------------
# 2) Near the end of the document: remove the duplicate bold title
#    paragraph entirely, and rewrite the italic paragraph's text
#    from the old meta description into the new "Prompt: ..." text,
#    preserving its italic run formatting.
# ------------------------------------------------------------------

$oldTitlePara = $d.Paragraphs.Item(49)
$oldTitlePara.Range.Delete()

$promptPara = $d.Paragraphs.Item(49)
$promptText = 'Prompt: Create a cartoon-style feature image for Novomatic''s "African Simba" slot game featuring a happy Maya warrior with glasses. The image should be colorful and eye-catching to represent the vibrant African savanna theme of the game. The Maya warrior could be holding a tribal spear or shield, and should be surrounded by the iconic animals of the savanna, such as a lion, elephant, giraffe, or buffalo. Make sure to include the game''s title in bold, African-inspired letters to tie in with the theme.'

# NOTE: Find.Execute's replacement text goes through Word's "smart quotes"
# AutoCorrect, turning straight quotes/apostrophes into curly ones. The
# target text uses straight quotes, so assign Range.Text directly instead
# (on the paragraph range minus its trailing paragraph mark) -- this keeps
# the run's existing formatting (italics) since it is a single run.
$promptRange = $promptPara.Range
$promptRange = $d.Range($promptRange.Start, $promptRange.End - 1)
$promptRange.Text = $promptText

Write-Output "done"
